$wb = $excel.ActiveWorkbook

$sheetData = @{}

$sheetData["ALC"] = @{
    19 = @{ "H" = 4546073; "I" = 6993341; "J" = 1146.5714; "K" = 6993341; "L" = 1146.5714; "M" = -6993166; "N" = -1496.5714 }
    62 = @{ "H" = 4555; "I" = 3362.7273; "J" = 6428.5713; "K" = 3362.7273; "L" = 6428.5713; "M" = -2738.7273; "N" = -7676.5713 }
    65 = @{ "H" = 4555; "I" = 3362.7273; "J" = 6428.5713; "K" = 16813.6365; "L" = 32142.8565; "M" = -13693.6365; "N" = -38382.85649999999 }
    86 = @{ "H" = 1840.2; "I" = 1800.25; "J" = 2000; "K" = 1800.25; "L" = 2000; "M" = -677.25; "N" = -4246 }
    88 = @{ "H" = 15301.333; "I" = 0; "J" = 15301.333; "K" = 0; "L" = 15301.333; "M" = $null; "N" = -16113.333 }
    89 = @{ "H" = 1840.2; "I" = 1800.25; "J" = 2000; "K" = 9001.25; "L" = 10000; "M" = -3385.25; "N" = -21232 }
    91 = @{ "H" = 15301.333; "I" = 0; "J" = 15301.333; "K" = 0; "L" = 15301.333; "M" = $null; "N" = -18109.333 }
    98 = @{ "H" = 2912.4707; "I" = 2557.1428; "J" = 4570.6665; "K" = 2557.1428; "L" = 4570.6665; "M" = -1059.1428; "N" = -7566.6665 }
    107 = @{ "H" = 808; "I" = 906.5; "J" = 718.4545; "K" = 906.5; "L" = 718.4545; "M" = 1013.5; "N" = -4558.4545 }
    113 = @{ "H" = 2125; "I" = 2000; "J" = 2138.889; "K" = 2000; "L" = 2138.889; "M" = 1254; "N" = -8646.889 }
    122 = @{ "H" = 2912.4707; "I" = 2557.1428; "J" = 4570.6665; "K" = 7671.428400000001; "L" = 13711.9995; "M" = -5221.428400000001; "N" = -18611.9995 }
    131 = @{ "H" = 2686.973; "I" = 1974.6; "J" = 5740; "K" = 5923.799999999999; "L" = 17220; "M" = -883.7999999999993; "N" = -27300 }
    135 = @{ "H" = 950.6667; "I" = 950.6667; "J" = 0; "K" = 8556.0003; "L" = 0; "M" = $null; "N" = -6021.0003 }
    137 = @{ "H" = 14316208; "I" = 25001030; "J" = 69780.4; "K" = 75003090; "L" = 209341.2; "M" = -75000540; "N" = -214441.2 }
}

$sheetData["ARM"] = @{
    97 = @{ "H" = 990.3077; "J" = 2499.6667; "L" = 2499.6667; "N" = -3491.6667 }
    122 = @{ "H" = 759.9375; "I" = 588.5714; "J" = 1959.5; "K" = 1765.7142; "L" = 5878.5; "M" = 684.2857999999999; "N" = -10778.5 }
}

$sheetData["BSM"] = @{
    94 = @{ "H" = 452.12195; "I" = 402.14285; "K" = 402.14285; "M" = 48.85714999999999 }
    99 = @{ "H" = 2302.68; "I" = 1529.6666; "J" = 3462.2; "K" = 1529.6666; "L" = 3462.2; "M" = -31.66660000000002; "N" = -6458.2 }
}

$sheetData["CRP"] = @{
    31 = @{ "H" = 2609.4238; "I" = 1762.9048; "J" = 4700.8237; "K" = 1762.9048; "L" = 4700.8237; "M" = -1467.9048; "N" = -5290.8237 }
    34 = @{ "H" = 2609.4238; "I" = 1762.9048; "J" = 4700.8237; "K" = 1762.9048; "L" = 4700.8237; "M" = -1560.9048; "N" = -5104.8237 }
}

$sheetData["CUL"] = @{
    3 = @{ "H" = 4823.3335; "J" = 5960; "L" = 17880; "N" = -18104 }
    113 = @{ "H" = 1149834.4; "I" = 2873956.2; "J" = 419.8889; "K" = 8621868.600000001; "L" = 1259.6667; "M" = -8619698.600000001; "N" = -5599.6667 }
    133 = @{ "H" = 3159; "I" = 2448.75; "K" = 7346.25; "M" = -2286.25 }
    138 = @{ "H" = 2169.1; "I" = 1620; "J" = 2266; "K" = 4860; "L" = 6798; "M" = 280; "N" = -17078 }
}

$sheetData["GSM"] = @{
    22 = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null; "N" = $null }
    113 = @{ "H" = 9927.652; "I" = 1708.909; "J" = 17461.5; "K" = 1708.909; "L" = 17461.5; "M" = 461.0909999999999; "N" = -21801.5 }
}

$sheetData["LTW"] = @{
    132 = @{ "H" = 1575.99; "I" = 1592.2812; "J" = 1185; "K" = 4776.8436; "L" = 3555; "M" = -2246.8436; "N" = -8615 }
}

$sheetData["WVR"] = @{
    132 = @{ "H" = 3178.6191; "I" = 3185.238; "J" = 3172; "K" = 9555.714; "L" = 9516; "M" = -7025.714; "N" = -14576 }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetData[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $address = "$colLetter$rowNum"
            $value = $cols[$colLetter]
            if ($null -eq $value) {
                $ws.Range($address).ClearContents()
            } else {
                $ws.Range($address).Value = $value
            }
        }
    }
}

Write-Host "Applied Shinryu_Profits updates to $($sheetData.Keys.Count) sheets"